$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Sales vs PO" (existing Sheet1, renamed)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Original (pre-edit) column A "ds" values, one per data row (rows 2..36).
# These become the new column C ("Order Week") values.
$origDs = @(45411,45418,45425,45432,45439,45446,45453,45460,45467,45474,45481,45488,45495,45502,45509,45516,45523,45530,45537,45544,45551,45558,45565,45572,45579,45586,45593,45600,45607,45614,45621,45628,45635,45642,45649)

# New column A "ds" values = original + 6 (one week later).
$newDs = @(45417,45424,45431,45438,45445,45452,45459,45466,45473,45480,45487,45494,45501,45508,45515,45522,45529,45536,45543,45550,45557,45564,45571,45578,45585,45592,45599,45606,45613,45620,45627,45634,45641,45648,45655)

# Insert a new column at C (shifts old PO_Requested_Qty column C -> D).
$ws1.Columns.Item(3).Insert()

# Header row.
$ws1.Cells.Item(1,3).Value = "Order Week"

# Fill column C (Order Week) with the original ds dates, using the same
# date style as column A. Copy format from A2 (date style) down column C.
$ws1.Range("A2").Copy()
$ws1.Range("C2:C36").PasteSpecial(-4122)

for ($i = 0; $i -lt $origDs.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 3).Value = $origDs[$i]
    $ws1.Cells.Item($row, 1).Value = $newDs[$i]
    $ws1.Cells.Item($row, 4).Value = 0
}

$ws1.Name = "Sales vs PO"

# ---------------------------------------------------------------------------
# Sheet 2: "Weekly Growth" — weeks that actually had a PO request, plus the
# week-over-week growth percentage.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"

$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"

# Copy the bold/bordered header style from sheet1's header row.
$ws1.Range("A1:C1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)

$growthDs  = @(45418,45425,45439,45446,45453,45460,45523,45544,45551,45558,45565,45572,45579,45586)
$growthQty = @(144,32,64,144,192,224,16,64,608,16,16,256,16,592)
$growthPct = @(0,-77.77777777777779,100,125,33.33333333333333,16.66666666666667,-92.85714285714286,300,850,-97.36842105263158,0,1500,-93.75,3600)

# Copy the date style down column A.
$ws1.Range("A2").Copy()
$ws2.Range("A2:A15").PasteSpecial(-4122)

for ($i = 0; $i -lt $growthDs.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 1).Value = $growthDs[$i]
    $ws2.Cells.Item($row, 2).Value = $growthQty[$i]
    $ws2.Cells.Item($row, 3).Value = $growthPct[$i]
}

# ---------------------------------------------------------------------------
# Sheet 3: "Volume Insights" — summary stats over the PO quantities.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"

$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"

$ws1.Range("A1:D1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

$ws3.Cells.Item(2,1).Value = 2384
$ws3.Cells.Item(2,2).Value = 170.2857142857143
$ws3.Cells.Item(2,3).Value = 608
$ws3.Cells.Item(2,4).Value = 16

# ---------------------------------------------------------------------------
# Sheet 4: "Prediction Info" — next week's forecast PO quantity.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"

$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws1.Range("A1").Copy()
$ws4.Range("A1").PasteSpecial(-4122)

$ws4.Cells.Item(2,1).Value = 282.3736263736264

